$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(188).Insert()

$ws.Cells.Item(188, 1).Value = 3
$ws.Cells.Item(188, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(188, 3).Value = "Coquimbo"
$ws.Cells.Item(188, 4).Value = 44461
$ws.Cells.Item(188, 5).Value = 5
$ws.Cells.Item(188, 6).Value = 100112040
$ws.Cells.Item(188, 7).Value = "Cilantro"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 230
$ws.Cells.Item(188, 11).Value = 2300
$ws.Cells.Item(188, 12).Value = 2500
$ws.Cells.Item(188, 13).Value = 2404
$ws.Cells.Item(188, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(188, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(188, 16).Value = 801
$ws.Cells.Item(188, 17).Value = 3
$ws.Cells.Item(188, 18).Value = "Hortaliza"
